$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 5
    3  = 2
    4  = 2
    5  = 3
    6  = 3
    7  = 3
    8  = 3
    9  = 3
    10 = 6
    11 = 3
    12 = 5
    13 = 1
    14 = 3
    15 = 1
    16 = 7
    17 = 7
    18 = 4
    19 = 5
    20 = 5
    21 = 7
    22 = 6
    23 = 6
    24 = 4
    25 = 6
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
